# QUERY_1000set_100times.xlsx — add INDEX() wrapper around the QUERY(...) calls
# used by the performance test, per commit "performance tests for QUERY
# function are added".
#
# The sheet holds 100 rows (A1:A100) of the same volatile formula:
#   =QUERY("ShuffledDataSet1000","ShuffledDataSet1000_100_local")
# stored as three formula groups (A1 alone, a shared group A2:A65, and a
# shared group A66:A100). All three get wrapped in INDEX(...,1,1) so the
# engine only needs to materialise the first element of the result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFormula = '=INDEX(QUERY("ShuffledDataSet1000","ShuffledDataSet1000_100_local"),1,1)'

# Re-write the formula across the whole used range in one shot so the
# engine keeps it as shared formulas instead of 100 distinct ones.
$ws.Range("A1:A100").Formula = $newFormula

# The longer formula text (this sheet shows formulas, not values) needs a
# wider column to keep fitting it ("bestFit"); widen column A to match.
$ws.Columns("A:A").ColumnWidth = 35

# Cursor ends up on A16 after the edit.
$ws.Range("A16").Select()
